# Cookie cutter worked on further, All dicts placed into one, more dicts added
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. BemptonDriveOdd: move selection from B21 to A9
# ---------------------------------------------------------------------------
$bemptonOdd = $wb.Worksheets.Item("BemptonDriveOdd")
$bemptonOdd.Range("A9").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. LynmouthDriveEven: keep its own selection (C7) - just make sure it is no
#    longer the tab that stays "active" at the end of the session; that is
#    handled naturally below because later sheets get activated afterwards.
# ---------------------------------------------------------------------------
$lynmouthEven = $wb.Worksheets.Item("LynmouthDriveEven")
$lynmouthEven.Range("C7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add the new "BeverleyRoadEven" sheet (even house numbers 62-136) after
#    the last existing sheet (BemptonDriveEven).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$beverleyEven = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$beverleyEven.Name = "BeverleyRoadEven"

$evenAddresses = @(
    "62 Beverley Rd Ruislip HA4 9AS UK",
    "64 Beverley Rd Ruislip HA4 9AS UK",
    "66 Beverley Rd Ruislip HA4 9AS UK",
    "68 Beverley Rd Ruislip HA4 9AS UK",
    "70 Beverley Rd Ruislip HA4 9AS UK",
    "72 Beverley Rd Ruislip HA4 9AS UK",
    "74 Beverley Rd Ruislip HA4 9AP UK",
    "76 Beverley Rd Ruislip HA4 9AS UK",
    "78 Beverley Rd Ruislip HA4 9AP UK",
    "80 Beverley Rd Ruislip HA4 9AP UK",
    "82 Beverley Rd Ruislip HA4 9AS UK",
    "84 Beverley Rd Ruislip HA4 9AP UK",
    "86 Beverley Rd Ruislip HA4 9AP UK",
    "88 Beverley Rd Ruislip HA4 9AS UK",
    "90 Beverley Rd Ruislip HA4 9AS UK",
    "92 Beverley Rd Ruislip HA4 9AS UK",
    "94 Beverley Rd Ruislip HA4 9AS UK",
    "96 Beverley Rd Ruislip HA4 9AS UK",
    "98 Beverley Rd Ruislip HA4 9AS UK",
    "100 Beverley Rd Ruislip HA4 9AS UK",
    "102 Beverley Rd Ruislip HA4 9AS UK",
    "104 Beverley Rd Ruislip HA4 9AS UK",
    "106 Beverley Rd Ruislip HA4 9AS UK",
    "108 Beverley Rd Ruislip HA4 9AS UK",
    "110 Beverley Rd Ruislip HA4 9AS UK",
    "112 Beverley Rd Ruislip HA4 9AS UK",
    "114 Beverley Rd Ruislip HA4 9AS UK",
    "116 Beverley Rd Ruislip HA4 9AS UK",
    "118 Beverley Rd Ruislip HA4 9AS UK",
    "120 Beverley Rd Ruislip HA4 9AR UK",
    "122 Beverley Rd Ruislip HA4 9AR UK",
    "124 Beverley Rd Ruislip HA4 9AR UK",
    "126 Beverley Rd Ruislip HA4 9AR UK",
    "128 Beverley Rd Ruislip HA4 9AS UK",
    "130 Beverley Rd Ruislip HA4 9AS UK",
    "132 Beverley Rd Ruislip HA4 9AS UK",
    "134 Beverley Rd Ruislip HA4 9AS UK",
    "136 Beverley Rd Ruislip HA4 9AS UK"
)

for ($i = 0; $i -lt $evenAddresses.Length; $i++) {
    $beverleyEven.Cells.Item($i + 1, 1).Value = $evenAddresses[$i]
}

$beverleyEven.Columns.Item(1).ColumnWidth = 37.33
$beverleyEven.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Add the new "BeverleyRoadOdd" sheet (odd house numbers 127-219) after
#    BeverleyRoadEven. This ends up being the final active sheet/tab.
# ---------------------------------------------------------------------------
$beverleyOdd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $beverleyEven)
$beverleyOdd.Name = "BeverleyRoadOdd"

# NOTE: these were originally typed in this (not numerically sorted) order,
# and then the column was sorted ascending - reproduce that same workflow so
# that the shared-string table ends up populated in the same order too.
$oddAddressesEntryOrder = @(
    "211 Beverley Rd Ruislip HA4 9AR UK",
    "145 Beverley Rd Ruislip HA4 9AP UK",
    "137 Beverley Rd Ruislip HA4 9AP UK",
    "135 Beverley Rd Ruislip HA4 9AP UK",
    "127 Beverley Rd Ruislip HA4 9AP UK",
    "129 Beverley Rd Ruislip HA4 9AP UK",
    "131 Beverley Rd Ruislip HA4 9AP UK",
    "133 Beverley Rd Ruislip HA4 9AP UK",
    "139 Beverley Rd Ruislip HA4 9AP UK",
    "141 Beverley Rd Ruislip HA4 9AP UK",
    "143 Beverley Rd Ruislip HA4 9AP UK",
    "147 Beverley Rd Ruislip HA4 9AP UK",
    "149 Beverley Rd Ruislip HA4 9AP UK",
    "151 Beverley Rd Ruislip HA4 9AP UK",
    "153 Beverley Rd Ruislip HA4 9AP UK",
    "155 Beverley Rd Ruislip HA4 9AP UK",
    "157 Beverley Rd Ruislip HA4 9AP UK",
    "159 Beverley Rd Ruislip HA4 9AP UK",
    "161 Beverley Rd Ruislip HA4 9AP UK",
    "163 Beverley Rd Ruislip HA4 9AP UK",
    "165 Beverley Rd Ruislip HA4 9AP UK",
    "167 Beverley Rd Ruislip HA4 9AP UK",
    "169 Beverley Rd Ruislip HA4 9AP UK",
    "171 Beverley Rd Ruislip HA4 9AS UK",
    "173 Beverley Rd Ruislip HA4 9AS UK",
    "175 Beverley Rd Ruislip HA4 9AP UK",
    "177 Beverley Rd Ruislip HA4 9AP UK",
    "219 Beverley Rd Ruislip HA4 9DT UK",
    "217 Beverley Rd Ruislip HA4 9AR UK",
    "215 Beverley Rd Ruislip HA4 9AR UK",
    "213 Beverley Rd Ruislip HA4 9AR UK",
    "209 Beverley Rd Ruislip HA4 9AR UK",
    "207 Beverley Rd Ruislip HA4 9AR UK",
    "205 Beverley Rd Ruislip HA4 9AR UK",
    "203 Beverley Rd Ruislip HA4 9AR UK",
    "201 Beverley Rd Ruislip HA4 9AR UK",
    "199 Beverley Rd Ruislip HA4 9AR UK",
    "197 Beverley Rd Ruislip HA4 9AR UK",
    "195 Beverley Rd Ruislip HA4 9AR UK",
    "193 Beverley Rd Ruislip HA4 9AR UK",
    "191 Beverley Rd Ruislip HA4 9AR UK",
    "189 Beverley Rd Ruislip HA4 9AR UK",
    "187 Beverley Rd Ruislip HA4 9AR UK",
    "185 Beverley Rd Ruislip HA4 9AR UK",
    "183 Beverley Rd Ruislip HA4 9AR UK",
    "181 Beverley Rd Ruislip HA4 9AR UK",
    "179 Beverley Rd Ruislip HA4 9AR UK"
)

for ($i = 0; $i -lt $oddAddressesEntryOrder.Length; $i++) {
    $beverleyOdd.Cells.Item($i + 1, 1).Value = $oddAddressesEntryOrder[$i]
}

# Sort the column ascending (text/lexicographic sort - matches original workflow)
$sortRange = $beverleyOdd.Range("A1:A47")
$sortRange.Sort($beverleyOdd.Range("A1")) | Out-Null

$beverleyOdd.Columns.Item(1).ColumnWidth = 34.17

# Scroll the view down a little and land the selection/active cell on A26,
# then select A26 last so BeverleyRoadOdd ends up as the active/selected tab.
$beverleyOdd.Activate()
$excel.ActiveWindow.ScrollRow = 5
$beverleyOdd.Range("A26").Select() | Out-Null
